$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight the reviewed requirement rows (11-19, columns A:B) in yellow ---
# Seed the two new cell styles in the exact order the target workbook uses:
#   style index 12 -> red font (fontId 3) + yellow fill
#   style index 13 -> normal font (fontId 2) + yellow fill
$ws.Range("A12").Interior.Color = 65535
$ws.Range("A11").Interior.Color = 65535
# Now paint the whole reviewed block; this reuses the two styles just created.
$ws.Range("A11:B19").Interior.Color = 65535

# --- Remove the stray blank row (old row 20: empty A cell, "Jan" in B) ---
$ws.Rows(20).Delete()

# --- Update the note on the header comment ---
$c = $ws.Range("D1").Comment
$c.Text($c.Text() + "`nKeltaisella pohjalla olevat vaatimukset on käyty läpi.")

# --- Restore the last-used selection ---
$ws.Range("J16").Select()
